# Updated cryptos list on Sat Dec  2 01:00:28 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (D) column cells are stored as text (quote-prefixed numbers like "38.689.00"
# or values with significant trailing zeros like "2.40"); force text format first so
# Excel does not silently convert the assigned string into a numeric value and drop
# formatting such as trailing zeros.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "38.743.48"
$ws.Range("E2").Value = "  +2.67%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.092.42"
$ws.Range("E3").Value = "  +2.04%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "228.32"
$ws.Range("E5").Value = "  +0.41%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.615"
$ws.Range("E6").Value = "  +1.04%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "60.32"
$ws.Range("E7").Value = "  +2.02%  "
$ws.Range("E9").Value = "  +2.43%  "
$ws.Range("E10").Value = "  +1.50%  "
$ws.Range("E11").Value = "  -0.20%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.401.86"
$ws.Range("E12").Value = "  +2.11%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "15.11"
$ws.Range("E13").Value = "  +5.23%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.91"
$ws.Range("E14").Value = "  +3.16%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.804"
$ws.Range("E15").Value = "  +5.90%  "
$ws.Range("E16").Value = "  +0.49%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.095.44"
$ws.Range("E17").Value = "  +2.27%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "38.673.35"
$ws.Range("E18").Value = "  +2.79%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "71.52"
$ws.Range("E19").Value = "  +3.19%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.04"
$ws.Range("E20").Value = "  +0.95%  "
$ws.Range("E21").Value = "  +1.10%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "227.13"
$ws.Range("E22").Value = "  +2.14%  "
$ws.Range("E23").Value = "  -0.47%  "
$ws.Range("E24").Value = "  -1.80%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.33"
$ws.Range("E25").Value = "  +2.48%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "171.33"
$ws.Range("E26").Value = "  +1.47%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.56"
$ws.Range("E27").Value = "  +3.10%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.139"
$ws.Range("E28").Value = "  +8.07%  "
$ws.Range("E29").Value = "  +13.96%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.34"
$ws.Range("E30").Value = "  +3.17%  "
$ws.Range("E31").Value = "  +1.11%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.40"
$ws.Range("E32").Value = "  +0.13%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.52"
$ws.Range("E33").Value = "  +3.75%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.70"
$ws.Range("E34").Value = "  +3.50%  "
$ws.Range("E35").Value = "  +1.55%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.50"
$ws.Range("E36").Value = "  +1.83%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.40"
$ws.Range("E37").Value = "  +1.83%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.56"
$ws.Range("E38").Value = "  +3.75%  "
$ws.Range("E39").Value = "  -0.21%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.00"
$ws.Range("E40").Value = "  -1.00%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0230"
$ws.Range("E41").Value = "  +7.17%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.544.41"
$ws.Range("E42").Value = "  +0.67%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "100.30"
$ws.Range("E43").Value = "  +1.35%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0924"
$ws.Range("E44").Value = "  +3.17%  "
$ws.Range("E45").Value = "  -0.87%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.71"
$ws.Range("E46").Value = "  +8.74%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.13"
$ws.Range("E47").Value = "  +1.87%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.10"
$ws.Range("E48").Value = "  +0.08%  "
$ws.Range("E49").Value = "  +3.14%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.97"
$ws.Range("E50").Value = "  +0.40%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.290.47"
$ws.Range("E51").Value = "  +2.16%  "
